$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6: ex 1.7 and material
$ws.Cells.Item(6, 1).Value = 210909
$ws.Cells.Item(6, 2).Value = 60
$ws.Cells.Item(6, 3).Value = 1
$ws.Cells.Item(6, 4).Value = "ex 1.7 and material"

# Row 7: ex 1.8
$ws.Cells.Item(7, 1).Value = 210910
$ws.Cells.Item(7, 2).Value = 20
$ws.Cells.Item(7, 3).Value = 1
$ws.Cells.Item(7, 4).Value = "ex 1.8"

# Update selection to match diff (active cell E7)
$ws.Range("E7").Select()

$wb.Save()
